# "Verify Get Users.xlsx" - add new test case (Invalid API key) + latest report values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix punctuation on the existing "Missing API key" error message.
$ws.Range("E2").Value = "Missing API key."

# 2. Insert a new row for the new "Invalid API key" test case above the old row 3,
#    shifting the existing test cases down by one.
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "aaaa"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 401
$ws.Range("D3").Clear()
$ws.Range("E3").Value = "Invalid API key."
$ws.Range("F3").Value = "Check authorization"

# 3. The old row 4 (duplicate "missing key" test case) is now at row 5 - it has
#    been superseded by the new row 3 above, so remove it and shift the rest up.
$ws.Rows.Item(5).Delete()

# 4. Match the formatting used by the rest of the data rows (reuse the existing
#    cell styles instead of minting new ones).
$ws.Range("A4").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B3").PasteSpecial(-4122)

# 5. The D:E merged column width is now split into two separate column widths.
$ws.Range("D1").ColumnWidth = 10.8
$ws.Range("E1").ColumnWidth = 35.67
